{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" copyright line, and the blank paragraph that\n// separates them from the preceding \"LOM3099: Est\u00e1tica (Requisito fraco)\"\n// paragraph \u2014 i.e. collapse the footer block back down to just the\n// requirements line, leaving the existing blank/page-break paragraphs after\n// it untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOM3099: Est\u00e1tica (Requisito fraco)\" paragraph.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"LOM3099: Est\u00e1tica (Requisito fraco)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The three paragraphs immediately after it are: a blank paragraph, the\n  // \"Ver no Jupiter...\" line, and the \"\u00a9 2020 ...\" line. Only delete when\n  // that exact pattern is found (so re-running this script against an\n  // already-edited document is a safe no-op).\n  const jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n  const copyrightText =\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n  const blankIdx = anchorIndex + 1;\n  const jupiterIdx = anchorIndex + 2;\n  const copyrightIdx = anchorIndex + 3;\n\n  const matches =\n    copyrightIdx < items.length &&\n    items[blankIdx].text.trim() === \"\" &&\n    items[jupiterIdx].text.trim() === jupiterText &&\n    items[copyrightIdx].text.trim() === copyrightText;\n\n  if (matches) {\n    items[copyrightIdx].delete();\n    items[jupiterIdx].delete();\n    items[blankIdx].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" copyright line, and the blank paragraph that\n# separates them from the preceding \"LOM3099: Est\u00e1tica (Requisito fraco)\"\n# paragraph \u2014 i.e. collapse the footer block back down to just the\n# requirements line, leaving the existing blank/page-break paragraphs after\n# it untouched.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOM3099: Est\u00e1tica (Requisito fraco)\"\n$blankText = \"\"\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -gt 0) {\n    $blankIdx = $anchorIndex + 1\n    $jupiterIdx = $anchorIndex + 2\n    $copyrightIdx = $anchorIndex + 3\n\n    # Only delete when the exact pattern is found (so re-running this script\n    # against an already-edited document is a safe no-op).\n    $patternMatches = $false\n    if ($copyrightIdx -le $count) {\n        $blankOk = $d.Paragraphs.Item($blankIdx).Range.Text.Trim() -eq $blankText\n        $jupiterOk = $d.Paragraphs.Item($jupiterIdx).Range.Text.Trim() -eq $jupiterText\n        $copyrightOk = $d.Paragraphs.Item($copyrightIdx).Range.Text.Trim() -eq $copyrightText\n        $patternMatches = $blankOk -and $jupiterOk -and $copyrightOk\n    }\n\n    if ($patternMatches) {\n        # Delete from the highest index down so earlier indices stay valid.\n        $d.Paragraphs.Item($copyrightIdx).Range.Delete()\n        $d.Paragraphs.Item($jupiterIdx).Range.Delete()\n        $d.Paragraphs.Item($blankIdx).Range.Delete()\n    }\n}\n"}
